$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Payment Date") for "Customer Name"
$ws.Range("D1").EntireColumn.Insert()

# New header and value
$ws.Range("D1").Value = "Customer Name"
$ws.Range("D2").Value = "KLIKLELANG-Eddy susiyanto"

# Match the style of the adjacent header/value cells
$ws.Range("D1").Style = $ws.Range("E1").Style
$ws.Range("D2").Style = $ws.Range("C2").Style

# New column width (same width as the VA Number column, no auto-fit)
$ws.Range("D1").EntireColumn.Width = $ws.Range("C1").EntireColumn.Width

# Update selection to mirror the target state
$ws.Range("G13").Select()
